$wb = $excel.ActiveWorkbook

# --- Generator sheet: p_nom_extendable False -> True, p_nom 100 -> 0 ---
$gen = $wb.Worksheets.Item("Generator")
$gen.Range("D2").Value = "'True"
$gen.Range("H2").Value = 0
$gen.Range("D3").Value = "'True"
$gen.Range("H3").Value = 0

# --- StorageUnit sheet: p_nom_extendable False -> True, p_nom 100 -> 0 ---
$su = $wb.Worksheets.Item("StorageUnit")
$su.Range("D2").Value = "'True"
$su.Range("J2").Value = 0

# --- Update selections / active sheet ---
$gen.Range("D4").Select()

$su.Activate()
$su.Range("D3").Select()
